$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("zh-cn")
$ws1.Range("E2").Value = "2016-03-23 16:23:06"
$ws1.Range("H2").Value = "2016-03-23 16:23:34"
$ws1.Range("E4").Value = "2016-03-23 16:23:06"
$ws1.Range("H4").Value = "2016-03-23 16:23:34"

$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("E2").Value = "2016-03-23 16:23:12"
$ws2.Range("H2").Value = "2016-03-23 16:23:41"
$ws2.Range("E4").Value = "2016-03-23 16:23:12"
$ws2.Range("H4").Value = "2016-03-23 16:23:41"
